$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Group = "tot") updated values
$ws.Range("B2").Value = -0.999490387984652
$ws.Range("C2").Value = -167965945.557761
$ws.Range("D2").Value = 0.956625068060372
$ws.Range("E2").Value = 0.999490018715337
$ws.Range("F2").Value = -0.956629251004605
$ws.Range("G2").Value = 123400.043840362
$ws.Range("H2").Value = 7053953.60237105
$ws.Range("I2").Value = -2938364.66570834
$ws.Range("J2").Value = 0.958218396245853
$ws.Range("K2").Value = 0.999569298810958
$ws.Range("L2").Value = -0.956629251004605
$ws.Range("M2").Value = 51.2251243362736
$ws.Range("N2").Value = 2923.55996441976
$ws.Range("O2").Value = -2938364.66570834

# Row 4 (Group = "w") updated values
$ws.Range("B4").Value = -0.999490392340614
$ws.Range("C4").Value = -167966143.768117
$ws.Range("D4").Value = 0.956629474753855
$ws.Range("E4").Value = 0.999490393066385
$ws.Range("F4").Value = -0.956629315659008
$ws.Range("G4").Value = 123411.930720036
$ws.Range("H4").Value = 7054603.24085416
$ws.Range("I4").Value = -2938368.31894859
$ws.Range("J4").Value = 0.956661154052383
$ws.Range("K4").Value = 0.999490528013075
$ws.Range("L4").Value = -0.956629315659008
$ws.Range("M4").Value = 50.2858913260332
$ws.Range("N4").Value = 2874.40047728179
$ws.Range("O4").Value = -2938368.31894859
